# Update Leve profit/price figures across all sheets to match refreshed market data
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 13197.4
$ws.Range("J32").Value = 13197.4
$ws.Range("L32").Value = 13197.4
$ws.Range("N32").Value = -13849.4
$ws.Range("H51").Value = 88611.55499999999
$ws.Range("I51").Value = 7998.5
$ws.Range("J51").Value = 111643.86
$ws.Range("K51").Value = 7998.5
$ws.Range("L51").Value = 111643.86
$ws.Range("M51").Value = -7514.5
$ws.Range("N51").Value = -112611.86
$ws.Range("H58").Value = 9037.666999999999
$ws.Range("J58").Value = 18062.5
$ws.Range("L58").Value = 54187.5
$ws.Range("N58").Value = -54487.5
$ws.Range("H98").Value = 5833.7
$ws.Range("I98").Value = 5162.5713
$ws.Range("K98").Value = 5162.5713
$ws.Range("M98").Value = -3664.5713
$ws.Range("H100").Value = 2375.4167
$ws.Range("J100").Value = 6002.5
$ws.Range("L100").Value = 6002.5
$ws.Range("N100").Value = -7084.5
$ws.Range("H113").Value = 6383.2354
$ws.Range("I113").Value = 6159.4
$ws.Range("K113").Value = 6159.4
$ws.Range("M113").Value = -2905.4
$ws.Range("H122").Value = 5833.7
$ws.Range("I122").Value = 5162.5713
$ws.Range("K122").Value = 15487.7139
$ws.Range("M122").Value = -13037.7139
$ws.Range("H137").Value = 17731.133
$ws.Range("I137").Value = 6799.4
$ws.Range("K137").Value = 20398.2
$ws.Range("M137").Value = -17848.2
$ws.Range("H138").Value = 5140.921
$ws.Range("I138").Value = 3950
$ws.Range("J138").Value = 5207.0835
$ws.Range("K138").Value = 11850
$ws.Range("L138").Value = 15621.2505
$ws.Range("M138").Value = -6710
$ws.Range("N138").Value = -25901.2505

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 292.6
$ws.Range("I5").Value = 290.75
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 290.75
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -178.75
$ws.Range("N5").Value = -524
$ws.Range("H137").Value = 79950.75
$ws.Range("J137").Value = 79950.75
$ws.Range("L137").Value = 79950.75
$ws.Range("N137").Value = -90150.75

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 292.6
$ws.Range("I4").Value = 290.75
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 290.75
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -175.75
$ws.Range("N4").Value = -530
$ws.Range("H94").Value = 63873.918
$ws.Range("I94").Value = 36641.855
$ws.Range("K94").Value = 36641.855
$ws.Range("M94").Value = -36190.855
$ws.Range("H105").Value = 96845.95
$ws.Range("I105").Value = 1639.3529
$ws.Range("J105").Value = 501474
$ws.Range("K105").Value = 1639.3529
$ws.Range("L105").Value = 501474
$ws.Range("M105").Value = 107.6470999999999
$ws.Range("N105").Value = -504968
$ws.Range("H107").Value = 2654.875
$ws.Range("I107").Value = 2675.2856
$ws.Range("K107").Value = 2675.2856
$ws.Range("M107").Value = -755.2856000000002

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 149999
$ws.Range("J9").Value = 149999
$ws.Range("L9").Value = 149999
$ws.Range("N9").Value = -150335
$ws.Range("H16").Value = 853.9524
$ws.Range("I16").Value = 853.9524
$ws.Range("K16").Value = 853.9524
$ws.Range("M16").Value = -566.9524
$ws.Range("H22").Value = 421.55554
$ws.Range("I22").Value = 485.2857
$ws.Range("J22").Value = 198.5
$ws.Range("K22").Value = 485.2857
$ws.Range("L22").Value = 198.5
$ws.Range("M22").Value = -135.2857
$ws.Range("N22").Value = -898.5
$ws.Range("H31").Value = 6496007.5
$ws.Range("I31").Value = 6995162
$ws.Range("J31").Value = 7000
$ws.Range("K31").Value = 6995162
$ws.Range("L31").Value = 7000
$ws.Range("M31").Value = -6994867
$ws.Range("N31").Value = -7590
$ws.Range("H34").Value = 6496007.5
$ws.Range("I34").Value = 6995162
$ws.Range("J34").Value = 7000
$ws.Range("K34").Value = 6995162
$ws.Range("L34").Value = 7000
$ws.Range("M34").Value = -6994960
$ws.Range("N34").Value = -7404
$ws.Range("H62").Value = 4020.6667
$ws.Range("I62").Value = 3931.6365
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 3931.6365
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -3307.6365
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 4020.6667
$ws.Range("I65").Value = 3931.6365
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 19658.1825
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -16538.1825
$ws.Range("N65").Value = -31240
$ws.Range("H69").Value = 13669.375
$ws.Range("I69").Value = 7510.75
$ws.Range("J69").Value = 19828
$ws.Range("K69").Value = 7510.75
$ws.Range("L69").Value = 19828
$ws.Range("M69").Value = -6761.75
$ws.Range("N69").Value = -21326
$ws.Range("H72").Value = 13669.375
$ws.Range("I72").Value = 7510.75
$ws.Range("J72").Value = 19828
$ws.Range("K72").Value = 22532.25
$ws.Range("L72").Value = 59484
$ws.Range("M72").Value = -18788.25
$ws.Range("N72").Value = -66972
$ws.Range("H97").Value = 66999
$ws.Range("J97").Value = 66999
$ws.Range("L97").Value = 66999
$ws.Range("N97").Value = -68981
$ws.Range("H99").Value = 4071.6
$ws.Range("I99").Value = 5257.5
$ws.Range("K99").Value = 5257.5
$ws.Range("M99").Value = -3759.5
$ws.Range("H113").Value = 853.9524
$ws.Range("I113").Value = 853.9524
$ws.Range("K113").Value = 853.9524
$ws.Range("M113").Value = 1316.0476
$ws.Range("H122").Value = 2360.75
$ws.Range("I122").Value = 2591.8
$ws.Range("K122").Value = 7775.400000000001
$ws.Range("M122").Value = -5325.400000000001
$ws.Range("H126").Value = 4071.6
$ws.Range("I126").Value = 5257.5
$ws.Range("K126").Value = 15772.5
$ws.Range("M126").Value = -13302.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6448.75
$ws.Range("I3").Value = 6448.75
$ws.Range("K3").Value = 19346.25
$ws.Range("M3").Value = -19234.25
$ws.Range("H56").Value = 6746.5
$ws.Range("I56").Value = 6746.5
$ws.Range("K56").Value = 6746.5
$ws.Range("M56").Value = -6216.5
$ws.Range("H64").Value = 23748.75
$ws.Range("I64").Value = 24997.5
$ws.Range("K64").Value = 74992.5
$ws.Range("M64").Value = -74722.5
$ws.Range("H67").Value = 23748.75
$ws.Range("I67").Value = 24997.5
$ws.Range("K67").Value = 74992.5
$ws.Range("M67").Value = -74056.5
$ws.Range("H131").Value = 1200706.2
$ws.Range("I131").Value = 168121.17
$ws.Range("J131").Value = 1820257.2
$ws.Range("K131").Value = 504363.51
$ws.Range("L131").Value = 5460771.6
$ws.Range("M131").Value = -499323.51
$ws.Range("N131").Value = -5470851.6
$ws.Range("H132").Value = 5589.737
$ws.Range("I132").Value = 2211.111
$ws.Range("J132").Value = 8630.5
$ws.Range("K132").Value = 19899.999
$ws.Range("L132").Value = 77674.5
$ws.Range("M132").Value = -17369.999
$ws.Range("N132").Value = -82734.5
$ws.Range("H133").Value = 4765.143
$ws.Range("I133").Value = 3392.8333
$ws.Range("K133").Value = 10178.4999
$ws.Range("M133").Value = -5118.499899999999
$ws.Range("H134").Value = 1685.4
$ws.Range("I134").Value = 1685.4
$ws.Range("K134").Value = 5056.200000000001
$ws.Range("M134").Value = 13.79999999999927
$ws.Range("H136").Value = 6328
$ws.Range("I136").Value = 1225
$ws.Range("J136").Value = 8596
$ws.Range("K136").Value = 3675
$ws.Range("L136").Value = 25788
$ws.Range("M136").Value = 1425
$ws.Range("N136").Value = -35988
$ws.Range("H137").Value = 7561.5835
$ws.Range("I137").Value = 4506.75
$ws.Range("K137").Value = 13520.25
$ws.Range("M137").Value = -8420.25
$ws.Range("H138").Value = 2964.3333
$ws.Range("I138").Value = 2964.3333
$ws.Range("K138").Value = 8892.999899999999
$ws.Range("M138").Value = -3752.999899999999
$ws.Range("H139").Value = 3487.6667
$ws.Range("I139").Value = 3373
$ws.Range("K139").Value = 10119
$ws.Range("M139").Value = -4979
$ws.Range("H140").Value = 1986
$ws.Range("I140").Value = 1837.7142
$ws.Range("K140").Value = 5513.142599999999
$ws.Range("M140").Value = -333.1425999999992

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 15724.5
$ws.Range("I41").Value = 15724.5
$ws.Range("K41").Value = 15724.5
$ws.Range("M41").Value = -15369.5
$ws.Range("H99").Value = 3999
$ws.Range("I99").Value = 3999
$ws.Range("K99").Value = 3999
$ws.Range("M99").Value = -1753

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 262.25
$ws.Range("I55").Value = 249.66667
$ws.Range("K55").Value = 249.66667
$ws.Range("M55").Value = -76.66667000000001
$ws.Range("H122").Value = 42542
$ws.Range("J122").Value = 41423.855
$ws.Range("L122").Value = 124271.565
$ws.Range("N122").Value = -129171.565
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2308.7693
$ws.Range("I122").Value = 2308.7693
$ws.Range("K122").Value = 6926.3079
$ws.Range("M122").Value = -4476.3079
$ws.Range("H126").Value = 2265.25
$ws.Range("I126").Value = 2774.5
$ws.Range("J126").Value = 1246.75
$ws.Range("K126").Value = 8323.5
$ws.Range("L126").Value = 3740.25
$ws.Range("M126").Value = -5853.5
$ws.Range("N126").Value = -8680.25
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("H136").Value = 1374.3
$ws.Range("I136").Value = 1427
$ws.Range("J136").Value = 900
$ws.Range("K136").Value = 4281
$ws.Range("L136").Value = 2700
$ws.Range("M136").Value = -1731
$ws.Range("N136").Value = -7800
